$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "C3" = "0.53"
    "E3" = "0.5"
    "H3" = "0.56"
    "J3" = "0.38"
    "C4" = "0.31"
    "E4" = "0.32"
    "E5" = "0.15"
    "F5" = "0.14"
    "J5" = "0.22"
    "E6" = "0.03"
    "C8" = "0.62"
    "D8" = "0.44"
    "E8" = "0.27"
    "F8" = "0.21"
    "H8" = "0.34"
    "J8" = "0.13"
    "C9" = "0.28"
    "D9" = "0.31"
    "E9" = "0.25"
    "F9" = "0.29"
    "H9" = "0.31"
    "J9" = "0.27"
    "C10" = "0.06"
    "D10" = "0.2"
    "E10" = "0.35"
    "F10" = "0.39"
    "H10" = "0.26"
    "J10" = "0.43"
    "C11" = "0.03"
    "D11" = "0.03"
    "E11" = "0.13"
    "F11" = "0.08"
    "H11" = "0.06"
    "J11" = "0.15"
    "C13" = "3.55"
    "D13" = "3.5"
    "E13" = "2.23"
    "F13" = "2.29"
    "H13" = "3.01"
    "J13" = "2.46"
    "D14" = "0.16"
    "E14" = "0.28"
    "J14" = "0.29"
    "D15" = "0.55"
    "E15" = "0.3"
    "F15" = "0.33"
    "J15" = "0.33"
    "D16" = "0.46"
    "E16" = "0.41"
    "F16" = "0.35"
    "H16" = "0.29"
    "C17" = "621535"
    "D17" = "16159"
    "E17" = "5267"
    "F17" = "6292"
    "H17" = "1955"
    "J17" = "6260"
    "C19" = "0.55"
    "E19" = "0.54"
    "H19" = "0.6"
    "J19" = "0.45"
    "C20" = "0.31"
    "E20" = "0.32"
    "F20" = "0.28"
    "H20" = "0.28"
    "J20" = "0.35"
    "E21" = "0.13"
    "H21" = "0.11"
    "C24" = "0.57"
    "D24" = "0.43"
    "E24" = "0.21"
    "F24" = "0.18"
    "H24" = "0.27"
    "J24" = "0.11"
    "C25" = "0.29"
    "D25" = "0.31"
    "E25" = "0.27"
    "F25" = "0.31"
    "H25" = "0.3"
    "J25" = "0.25"
    "C26" = "0.06"
    "D26" = "0.19"
    "E26" = "0.38"
    "F26" = "0.38"
    "H26" = "0.31"
    "J26" = "0.43"
    "C27" = "0.08"
    "D27" = "0.05"
    "E27" = "0.14"
    "F27" = "0.09"
    "H27" = "0.1"
    "J27" = "0.19"
    "C29" = "3.71"
    "D29" = "3.62"
    "E29" = "2.4"
    "F29" = "2.43"
    "H29" = "3.18"
    "J29" = "2.5"
    "E30" = "0.19"
    "H30" = "0.19"
    "C31" = "0.56"
    "E31" = "0.24"
    "F31" = "0.3"
    "J31" = "0.27"
    "C32" = "0.75"
    "D32" = "0.75"
    "E32" = "0.67"
    "F32" = "0.62"
    "H32" = "0.68"
    "J32" = "0.68"
    "C33" = "568861"
    "D33" = "13698"
    "E33" = "3989"
    "F33" = "4681"
    "H33" = "1161"
    "J33" = "4833"
}

foreach ($cellRef in $changes.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$cellRef]
}
